$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (D1-style "General, explicitly applied" - numFmtId 0 with
# applyNumberFormat flag - matches D1's cellXfs entry)
$ws.Range("E1").Value = "Fecha_inicio"
$ws.Range("F1").Value = "Fecha_fin"
$ws.Range("E1:F1").NumberFormat = "general"

# New data cells - dates stored as text (numFmtId 49 = "@")
$ws.Range("E2:F2").NumberFormat = "@"
$ws.Range("E2").Value = "2021-11-26"
$ws.Range("F2").Value = "2021-11-30"

# Update selection to match target state
$ws.Range("F3").Select()
